$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.784.37'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '1.770.77'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  +1.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4317'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3664'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.128'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07455'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.013'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.200'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.339'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = '1.770.22'
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001073'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06596'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.010'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.173'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("D23").Value = '27.851.96'
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.420'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.364'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("D29").Value = '1.977.43'
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.308'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.982'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.739'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09146'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2210'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.40%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06226'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6562'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.129'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02285'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.196'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.435'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.10%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.009'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.788'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5979'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.979'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06947'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.126'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.38%  '
